$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.436.11"
$ws.Range("E2").Value = "  -4.11%  "

$ws.Range("D3").Value = "2.351.87"
$ws.Range("E3").Value = "  -5.55%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "511.58"
$ws.Range("E5").Value = "  -4.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.04"
$ws.Range("E6").Value = "  -5.82%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.551"
$ws.Range("E8").Value = "  -2.43%  "

$ws.Range("D9").Value = "2.364.70"
$ws.Range("E9").Value = "  -5.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0955"
$ws.Range("E10").Value = "  -3.92%  "

$ws.Range("E11").Value = "  -0.95%  "

$ws.Range("E12").Value = "  -8.56%  "

$ws.Range("E13").Value = "  -5.34%  "

$ws.Range("D14").Value = "2.770.74"
$ws.Range("E14").Value = "  -5.82%  "

$ws.Range("D15").Value = "56.383.11"
$ws.Range("E15").Value = "  -4.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.37"
$ws.Range("E16").Value = "  -4.56%  "

$ws.Range("E17").Value = "  -4.41%  "

$ws.Range("D18").Value = "2.382.03"
$ws.Range("E18").Value = "  -5.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.25"
$ws.Range("E19").Value = "  -3.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.03"
$ws.Range("E20").Value = "  -4.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "309.72"
$ws.Range("E21").Value = "  -3.66%  "

$ws.Range("E22").Value = "  -0.70%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  +0.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.84"
$ws.Range("E24").Value = "  -1.53%  "

$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.389"
$ws.Range("E26").Value = "  -5.11%  "

$ws.Range("D27").Value = "2.465.62"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.153"
$ws.Range("E28").Value = "  -4.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.16"
$ws.Range("E29").Value = "  -4.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.57"
$ws.Range("E30").Value = "  +1.81%  "

$ws.Range("E31").Value = "  -4.44%  "

$ws.Range("D32").Value = "0.0₃0711"
$ws.Range("E32").Value = "  -6.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.11"
$ws.Range("E33").Value = "  -2.93%  "

$ws.Range("E34").Value = "  -7.00%  "

$ws.Range("E35").Value = "  -0.07%  "

$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.58"
$ws.Range("E37").Value = "  -3.34%  "

$ws.Range("E38").Value = "  -5.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.70"
$ws.Range("E39").Value = "  -6.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.800"
$ws.Range("E40").Value = "  +2.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.41"
$ws.Range("E41").Value = "  -3.48%  "

$ws.Range("E42").Value = "  -6.47%  "

$ws.Range("E43").Value = "  -3.94%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.82"
$ws.Range("E44").Value = "  -6.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "253.05"
$ws.Range("E45").Value = "  -9.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.566"
$ws.Range("E46").Value = "  -4.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "120.68"
$ws.Range("E47").Value = "  -8.45%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0904"
$ws.Range("E48").Value = "  -2.76%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0486"
$ws.Range("E49").Value = "  -4.48%  "

$ws.Range("E50").Value = "  -5.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.55"
$ws.Range("E51").Value = "  -6.77%  "
